# Add a new "appPrefix" config row above the existing config rows on Sheet1.
# This inserts a new row 4 (A4="appPrefix", B4="/mis_dashboard") and pushes
# the previously existing rows 4-26 down to rows 5-27, preserving their
# values, styles and hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift existing rows 4..26 down to 5..27 (bottom-up so we never
#        overwrite a row before it has been copied down). ---
for ($r = 26; $r -ge 4; $r--) {
    $dst = $r + 1
    $ws.Range("A$dst").Value2 = $ws.Range("A$r").Value2
    $ws.Range("B$dst").Value2 = $ws.Range("B$r").Value2
    $ws.Range("B$dst").Style = $ws.Range("B$r").Style.Name
}

# --- 2. Populate the newly freed row 4 with the new config entry. ---
$ws.Range("A4").Value2 = "appPrefix"
$ws.Range("B4").Value2 = "/mis_dashboard"
$ws.Range("B4").Style = "Normal"

# --- 3. Hyperlinks are anchored to the sheet, not to the cells that were
#        shifted via Value2 above, so rebuild them at the new, shifted
#        locations (each one row lower than before). ---
$ws.Cells.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B5"), "http://google.com")
$ws.Hyperlinks.Add($ws.Range("B6"), "http://google.com")
$ws.Hyperlinks.Add($ws.Range("B7"), "http://google.com")
$ws.Hyperlinks.Add($ws.Range("B8"), "http://google.com")
$ws.Hyperlinks.Add($ws.Range("B9"), "http://google.com")
$ws.Hyperlinks.Add($ws.Range("B10"), "http://google.com")
$ws.Hyperlinks.Add($ws.Range("B11"), "http://google.com")
$ws.Hyperlinks.Add($ws.Range("B12"), "http://google.com")
$ws.Hyperlinks.Add($ws.Range("B13:B16"), "http://google.com", "", "", "http://google.com")
$ws.Hyperlinks.Add($ws.Range("B17"), "http://google.com")
$ws.Hyperlinks.Add($ws.Range("B19"), "http://google.com")
$ws.Hyperlinks.Add($ws.Range("B20"), "http://google.com")
$ws.Hyperlinks.Add($ws.Range("B22"), "http://google.com")
$ws.Hyperlinks.Add($ws.Range("B23:B24"), "http://google.com", "", "", "http://google.com")
$ws.Hyperlinks.Add($ws.Range("B21"), "http://google.com")
$ws.Hyperlinks.Add($ws.Range("B27"), "http://google.com")

# Adding a hyperlink re-styles the cell with a freshly generated style, so
# re-apply the original "Hyperlink" cell style to keep the same style index
# used before the edit.
$ws.Range("B5").Style = "Hyperlink"
$ws.Range("B6").Style = "Hyperlink"
$ws.Range("B7").Style = "Hyperlink"
$ws.Range("B8").Style = "Hyperlink"
$ws.Range("B9").Style = "Hyperlink"
$ws.Range("B10").Style = "Hyperlink"
$ws.Range("B11").Style = "Hyperlink"
$ws.Range("B12").Style = "Hyperlink"
$ws.Range("B13:B16").Style = "Hyperlink"
$ws.Range("B17").Style = "Hyperlink"
$ws.Range("B19").Style = "Hyperlink"
$ws.Range("B20").Style = "Hyperlink"
$ws.Range("B22").Style = "Hyperlink"
$ws.Range("B23:B24").Style = "Hyperlink"
$ws.Range("B21").Style = "Hyperlink"
$ws.Range("B27").Style = "Hyperlink"

# --- 4. Match the new cell selection recorded in the saved workbook. ---
$ws.Range("G9").Select() | Out-Null
